$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "b"
$ws.Range("J2").Value = "Acknowledge (Backchannel)"
$ws.Range("I24").Value = "aa"
$ws.Range("J24").Value = "Agree/Accept"
$ws.Range("I30").Value = "sd"
$ws.Range("J30").Value = "Statement-non-opinion"
$ws.Range("I39").Value = "sd"
$ws.Range("J39").Value = "Statement-non-opinion"
$ws.Range("I55").Value = "b"
$ws.Range("J55").Value = "Acknowledge (Backchannel)"
$ws.Range("I65").Value = "aa"
$ws.Range("J65").Value = "Agree/Accept"
$ws.Range("I91").Value = "sd"
$ws.Range("J91").Value = "Statement-non-opinion"
$ws.Range("I102").Value = "aa"
$ws.Range("J102").Value = "Agree/Accept"
$ws.Range("I103").Value = "sd"
$ws.Range("J103").Value = "Statement-non-opinion"
$ws.Range("I111").Value = "aa"
$ws.Range("J111").Value = "Agree/Accept"
$ws.Range("I114").Value = "%"
$ws.Range("J114").Value = "Uninterpretable"
$ws.Range("I123").Value = "sv"
$ws.Range("J123").Value = "Statement-opinion"
$ws.Range("I128").Value = "sd"
$ws.Range("J128").Value = "Statement-non-opinion"
$ws.Range("I129").Value = "sd"
$ws.Range("J129").Value = "Statement-non-opinion"
$ws.Range("I144").Value = "b"
$ws.Range("J144").Value = "Acknowledge (Backchannel)"
$ws.Range("I153").Value = "aa"
$ws.Range("J153").Value = "Agree/Accept"
$ws.Range("I159").Value = "sd"
$ws.Range("J159").Value = "Statement-non-opinion"
$ws.Range("I162").Value = "b"
$ws.Range("J162").Value = "Acknowledge (Backchannel)"
$ws.Range("I164").Value = "b"
$ws.Range("J164").Value = "Acknowledge (Backchannel)"
$ws.Range("I167").Value = "sd"
$ws.Range("J167").Value = "Statement-non-opinion"
$ws.Range("I212").Value = "sv"
$ws.Range("J212").Value = "Statement-opinion"
$ws.Range("I216").Value = "%"
$ws.Range("J216").Value = "Uninterpretable"
$ws.Range("I221").Value = "b"
$ws.Range("J221").Value = "Acknowledge (Backchannel)"
$ws.Range("I222").Value = "sv"
$ws.Range("J222").Value = "Statement-opinion"
$ws.Range("I229").Value = "sd"
$ws.Range("J229").Value = "Statement-non-opinion"
$ws.Range("I234").Value = "aa"
$ws.Range("J234").Value = "Agree/Accept"
$ws.Range("I263").Value = "qy"
$ws.Range("J263").Value = "Yes-No-Question"
$ws.Range("I272").Value = "sd"
$ws.Range("J272").Value = "Statement-non-opinion"
$ws.Range("I277").Value = "sd"
$ws.Range("J277").Value = "Statement-non-opinion"
$ws.Range("I312").Value = "sv"
$ws.Range("J312").Value = "Statement-opinion"
$ws.Range("I313").Value = "sv"
$ws.Range("J313").Value = "Statement-opinion"
$ws.Range("I329").Value = "qy"
$ws.Range("J329").Value = "Yes-No-Question"
$ws.Range("I340").Value = "b"
$ws.Range("J340").Value = "Acknowledge (Backchannel)"
$ws.Range("I345").Value = "sd"
$ws.Range("J345").Value = "Statement-non-opinion"
$ws.Range("I351").Value = "aa"
$ws.Range("J351").Value = "Agree/Accept"
$ws.Range("I362").Value = "sv"
$ws.Range("J362").Value = "Statement-opinion"
$ws.Range("I369").Value = "sv"
$ws.Range("J369").Value = "Statement-opinion"
$ws.Range("I370").Value = "sv"
$ws.Range("J370").Value = "Statement-opinion"
$ws.Range("I380").Value = "b"
$ws.Range("J380").Value = "Acknowledge (Backchannel)"
$ws.Range("I388").Value = "b"
$ws.Range("J388").Value = "Acknowledge (Backchannel)"
$ws.Range("I395").Value = "ba"
$ws.Range("J395").Value = "Appreciation"
$ws.Range("I405").Value = "b"
$ws.Range("J405").Value = "Acknowledge (Backchannel)"
$ws.Range("I416").Value = "sd"
$ws.Range("J416").Value = "Statement-non-opinion"
$ws.Range("I419").Value = "ba"
$ws.Range("J419").Value = "Appreciation"
$ws.Range("I455").Value = "b"
$ws.Range("J455").Value = "Acknowledge (Backchannel)"
$ws.Range("I457").Value = "sd"
$ws.Range("J457").Value = "Statement-non-opinion"
$ws.Range("I463").Value = "sv"
$ws.Range("J463").Value = "Statement-opinion"
